$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: "Tienda Velázquez" -> "Tienda - Velázquez"
$ws.Name = "Tienda - Velázquez"

# Insert a new column before column B (shifts Nombre_TPV..total_operaciones right by one,
# bringing in a new "serie" column)
$ws.Range("B1").EntireColumn.Insert()

# New header for inserted column
$ws.Range("B1").Value = "serie"

# Protect the "fecha" (D) and "total_operaciones" (H) columns from being auto-converted
# to dates / numbers, since the source data stores every value as plain text.
$ws.Range("D2:D11").NumberFormat = "@"
$ws.Range("H2:H11").NumberFormat = "@"

# Update column A (Tienda) text on all data rows to the new store name
$ws.Range("A2:A11").Value = "Tienda - Velázquez"

# Fill in the new "serie" column values for each data row
$ws.Range("B2").Value = "V2"
$ws.Range("B3").Value = "V2"
$ws.Range("B4").Value = "V2"
$ws.Range("B5").Value = "V2"
$ws.Range("B6").Value = "V1"
$ws.Range("B7").Value = "V1"
$ws.Range("B8").Value = "V1"
$ws.Range("B9").Value = "V1"
$ws.Range("B10").Value = "V1"
$ws.Range("B11").Value = "V2"

# Row 2: BAR / 01/02/2025 / Mañana / EUROS / 72,75 / 9
$ws.Range("C2").Value = "BAR"
$ws.Range("D2").Value = "01/02/2025"
$ws.Range("E2").Value = "Mañana"
$ws.Range("F2").Value = "EUROS"
$ws.Range("G2").Value = "72,75"
$ws.Range("H2").Value = "9"

# Row 3: BAR / 01/02/2025 / Mañana / TARJETA VISA / 202,10 / 24
$ws.Range("C3").Value = "BAR"
$ws.Range("D3").Value = "01/02/2025"
$ws.Range("E3").Value = "Mañana"
$ws.Range("F3").Value = "TARJETA VISA"
$ws.Range("G3").Value = "202,10"
$ws.Range("H3").Value = "24"

# Row 4: BAR / 01/02/2025 / Mañana / EUROS / 83,40 / 9
$ws.Range("C4").Value = "BAR"
$ws.Range("D4").Value = "01/02/2025"
$ws.Range("E4").Value = "Mañana"
$ws.Range("F4").Value = "EUROS"
$ws.Range("G4").Value = "83,40"
$ws.Range("H4").Value = "9"

# Row 5: BAR / 01/02/2025 / Mañana / TARJETA VISA / 139,80 / 14
$ws.Range("C5").Value = "BAR"
$ws.Range("D5").Value = "01/02/2025"
$ws.Range("E5").Value = "Mañana"
$ws.Range("F5").Value = "TARJETA VISA"
$ws.Range("G5").Value = "139,80"
$ws.Range("H5").Value = "14"

# Row 6: SERVIDOR TIENDA / 01/02/2025 / Mañana / EUROS / 806,32 / 77
$ws.Range("C6").Value = "SERVIDOR TIENDA"
$ws.Range("D6").Value = "01/02/2025"
$ws.Range("E6").Value = "Mañana"
$ws.Range("F6").Value = "EUROS"
$ws.Range("G6").Value = "806,32"
$ws.Range("H6").Value = "77"

# Row 7: SERVIDOR TIENDA / 01/02/2025 / Mañana / TARJETA VISA / 2231,03 / 154
$ws.Range("C7").Value = "SERVIDOR TIENDA"
$ws.Range("D7").Value = "01/02/2025"
$ws.Range("E7").Value = "Mañana"
$ws.Range("F7").Value = "TARJETA VISA"
$ws.Range("G7").Value = "2231,03"
$ws.Range("H7").Value = "154"

# Row 8: SERVIDOR TIENDA / 01/02/2025 / Mañana / EUROS / 702,76 / 79
$ws.Range("C8").Value = "SERVIDOR TIENDA"
$ws.Range("D8").Value = "01/02/2025"
$ws.Range("E8").Value = "Mañana"
$ws.Range("F8").Value = "EUROS"
$ws.Range("G8").Value = "702,76"
$ws.Range("H8").Value = "79"

# Row 9: SERVIDOR TIENDA / 01/02/2025 / Mañana / SMS / 4,70 / 1
$ws.Range("C9").Value = "SERVIDOR TIENDA"
$ws.Range("D9").Value = "01/02/2025"
$ws.Range("E9").Value = "Mañana"
$ws.Range("F9").Value = "SMS"
$ws.Range("G9").Value = "4,70"
$ws.Range("H9").Value = "1"

# Row 10: SERVIDOR TIENDA / 01/02/2025 / Mañana / TARJETA VISA / 1809,80 / 154
$ws.Range("C10").Value = "SERVIDOR TIENDA"
$ws.Range("D10").Value = "01/02/2025"
$ws.Range("E10").Value = "Mañana"
$ws.Range("F10").Value = "TARJETA VISA"
$ws.Range("G10").Value = "1809,80"
$ws.Range("H10").Value = "154"

# Row 11: BAR / 01/02/2025 / Mañana / TARJETA VISA / 93,20 / 13
$ws.Range("C11").Value = "BAR"
$ws.Range("D11").Value = "01/02/2025"
$ws.Range("E11").Value = "Mañana"
$ws.Range("F11").Value = "TARJETA VISA"
$ws.Range("G11").Value = "93,20"
$ws.Range("H11").Value = "13"

Write-Host "edit complete"
